# Updates the "loading_percent" results sheet (Case_0_220) with the
# recalculated loading-percent values for the 380 kV case (rows 2-25,
# columns C,D,E,F,G,I,J,L). Column order per data row below matches
# ['C', 'D', 'E', 'F', 'G', 'I', 'J', 'L'].

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colIndices = @(3,4,5,6,7,9,10,12)

$data = @(
    @(2,5.2432754109253,9.63159552473418,9.745651487040339,98.45660804162037,3.982324374123778,73.95804897221228,12.81920940909195,10.7835427499928),
    @(3,5.104302578158464,9.507952370026249,9.769517484570104,97.91541962454778,3.995125036841528,73.60855874301228,12.86737087250441,10.82757635981412),
    @(4,5.019783746689234,9.435544475530445,9.784905445943942,97.61875603217129,4.003331750578287,73.42041106809553,12.89946642057406,10.85617237651303),
    @(5,4.985594426278229,9.406939693557337,9.791361987479481,97.50679875096158,4.006764203830502,73.35035957527111,12.91317891532743,10.86821925778727),
    @(6,4.979934027558397,9.402244896250163,9.79244534984664,97.48874706387835,4.007339509630406,73.33912641365818,12.91549407390725,10.8702434657943),
    @(7,5.019321568819845,9.435155025069399,9.784991767126817,97.61720999822525,4.003377683767833,73.41943956831911,12.89964878968761,10.8563332483378),
    @(8,5.195218212051413,9.588247358543912,9.753728849786581,98.26260144000041,3.986666547945113,73.83204242442808,12.83529069276932,10.79840295209612),
    @(9,5.544512110489363,9.915351719947514,9.698191484204274,99.81138231629728,3.956608911971034,74.85184375658845,12.72918208221883,10.69709614158885),
    @(10,5.801171664348588,10.17068383713276,9.660825974703121,101.1218036377911,3.936120989532283,75.73024894881925,12.66357806392296,10.6300530449684),
    @(11,5.917438816810465,10.28978353325452,9.644556442143589,101.7551583577401,3.92713383011899,76.15789633293811,12.63644007479573,10.6011340615302),
    @(12,5.961358400866619,10.3352799691126,9.638498886602877,102.0003087632881,3.923777404984851,76.32385627121393,12.62655492179062,10.59040842936013),
    @(13,5.951904915395127,10.32546433051991,9.639798914352074,101.9472757098008,3.924498206141678,76.28793538448919,12.62866641910775,10.59270839106606),
    @(14,5.921054482390093,10.29351877632686,9.644056019713517,101.7752208246551,3.926856762878582,76.17146951648057,12.63561895879296,10.60024715078333),
    @(15,5.902142534163008,10.27400194071082,9.646677039917797,101.6705227220402,3.928307512086429,76.10065366282352,12.63992863608518,10.60489415541516),
    @(16,5.793560172666507,10.16295732205418,9.661903765552532,101.0811583391498,3.936714931062332,75.70286458618061,12.66540623682013,10.63197459920527),
    @(17,5.726795134269773,10.09557019533218,9.671430429916111,100.7291127173016,3.941957135118663,75.46601090799787,12.68173051875504,10.64899082976388),
    @(18,5.688349722356572,10.05709005267103,9.676978542035917,100.530133473015,3.945003711349318,75.33242109138912,12.69137432942903,10.65892686745577),
    @(19,5.675326436906555,10.04411013268174,9.678868864969365,100.4633661729457,3.946040654073486,75.2876440313084,12.69468320196381,10.66231664001291),
    @(20,5.733907219366788,10.10271498144356,9.670409207343024,100.7662257780599,3.941395851380088,75.49095102481419,12.67996641312114,10.64716404023999),
    @(21,5.930119226100949,10.30289143570872,9.642802810367801,101.8256137056621,3.926162735560358,76.20556943868012,12.63356618701989,10.59802673204474),
    @(22,6.057708953831865,10.43601350788488,9.625362264286167,102.5489264674522,3.916479385065404,76.6960293919653,12.60552361654176,10.56722540899457),
    @(23,5.989682680815302,10.36476292712541,9.634615999907858,102.1600656514099,3.921623006817365,76.43212575397429,12.62028078560832,10.58354511615123),
    @(24,5.730692033712216,10.09948400665169,9.670870680487166,100.7494363186423,3.941649505628677,75.47966755916437,12.68076315946355,10.64798945453068),
    @(25,5.449832611648883,9.82412027549972,9.712606294651771,99.36204316489031,3.964455808104168,74.55329849327032,12.7557255898399,10.72319724458163)
)

foreach ($rowData in $data) {
    $rowIndex = $rowData[0]
    for ($i = 0; $i -lt $colIndices.Length; $i++) {
        $colIndex = $colIndices[$i]
        $value = $rowData[$i + 1]
        $ws.Cells.Item($rowIndex, $colIndex).Value = $value
    }
}

Write-Output "Updated $($data.Length) rows across $($colIndices.Length) columns."
